$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update G9 ("Test data" for TC_Vtiger_004): "Geeta Vishwas" -> "Narendra Modi" ---
$g9Text = @"
*URL:http://localhost:8888/
*User Name: admin                                                                                                                                                                                                                                                                                                                                                                                                               *Password: root    




Mr.
 Shin                                                                                                                                                                                                                                                                                                                                                                                         Chan
Mewtwo
Public Relations

Hero
Super Heros
shinchan@gmail.com
Narendra Modi
9999999999

8888888888
7777777777
6666666666
5555555555
Not Applicable
26-01-2000
God
Not Applicable
"@
$ws.Range("G9").Value = $g9Text

# --- Update H6 ("Expected result" for TC_Vtiger_001) with new rich text,
#     bolding "Creating New Organization " and "Organization Information " ---
$h6Text = @"
*Homepage should be displayed.          
*Creating New Organization page should display.                                 *Organization Information page should be displayed.
"@
$ws.Range("H6").Value = $h6Text
$ws.Range("H6").Characters(43, 26).Font.Bold = $true
$ws.Range("H6").Characters(69, 54).Font.Bold = $false
$ws.Range("H6").Characters(123, 25).Font.Bold = $true
$ws.Range("H6").Characters(148, 25).Font.Bold = $false

# --- Update selection to H6 ---
$ws.Range("H6").Select()
